# Fichero de Preguntas - Corrección
# Row 346 is a duplicate question ("Why is it important to decouple
# deployment from release?", already present earlier at row 200 with the
# same answer options just in a different order). Delete that whole row so
# every row below shifts up by one and the Nº numbering/answer table stays
# consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(346).Select()
$ws.Rows.Item(346).Delete()
